# Insert a new weekly price record for "Betarraga" (Feria Lagunitas de Puerto Montt)
# as row 243, pushing the existing rows 243-308 down to 244-309.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 243 (shifts old rows 243..308 -> 244..309)
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with the new data record
$ws.Range("A243").Value = 4
$ws.Range("B243").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C243").Value = "Los Lagos"
$ws.Range("D243").Value = 44736
$ws.Range("E243").Value = 10
$ws.Range("F243").Value = 100114014
$ws.Range("G243").Value = "Betarraga"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 1000
$ws.Range("K243").Value = 1000
$ws.Range("L243").Value = 1200
$ws.Range("M243").Value = 1100
$ws.Range("N243").Value = "$/paquete 5 unidades"
$ws.Range("O243").Value = "Región del Maule"
$ws.Range("P243").Value = 220
$ws.Range("Q243").Value = 5
$ws.Range("R243").Value = "Hortaliza"
